$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Thbs2", "Itgb1", "ECs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 295.3921834609324, 2658.529651148392, 0.003133449643399065, 0.003133449643399065),
    @("ECs", "Thbs2", "Itgb1", "FAPs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 324.7732161927742, 2922.958945734968, 0.003445116612570787, 0.003445116612570787),
    @("ECs", "Thbs2", "Itgb1", "M2", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 131.0068261449578, 1179.06143530462, 0.001389689083364161, 0.001389689083364161),
    @("ECs", "Thbs2", "Itgb1", "sCs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 180.0160174796526, 1620.144157316874, 0.001909566865205625, 0.001909566865205625),
    @("FAPs", "Thbs2", "Itgb1", "ECs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 28898.91150581132, 260090.2035523019, 0.3065527424982858, 0.3065527424982858),
    @("FAPs", "Thbs2", "Itgb1", "FAPs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 31773.32698600001, 285959.9428740001, 0.3370438545373735, 0.3370438545373734),
    @("FAPs", "Thbs2", "Itgb1", "M2", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 12816.70567942114, 115350.3511147902, 0.1359565489181078, 0.1359565489181078),
    @("FAPs", "Thbs2", "Itgb1", "sCs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 17611.39004364041, 158502.5103927637, 0.1868174140745515, 0.1868174140745515),
    @("M2", "Thbs2", "Itgb1", "ECs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 13.047049422004, 117.423444798036, 0.0001383999802560631, 0.0001383999802560631),
    @("M2", "Thbs2", "Itgb1", "FAPs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 14.344767532316, 129.102907790844, 0.0001521658636397967, 0.0001521658636397967),
    @("M2", "Thbs2", "Itgb1", "M2", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 5.786383767189998, 52.07745390471, 0.00006138057527263457, 0.00006138057527263458),
    @("M2", "Thbs2", "Itgb1", "sCs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 7.951049514213, 71.55944562791699, 0.00008434283186865066, 0.00008434283186865065),
    @("sCs", "Thbs2", "Itgb1", "ECs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 697.2352463801047, 6275.117217420943, 0.007396104759907931, 0.007396104759907931),
    @("sCs", "Thbs2", "Itgb1", "FAPs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 766.5853942265086, 6899.268548038577, 0.008131754544181002, 0.008131754544181),
    @("sCs", "Thbs2", "Itgb1", "M2", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 309.2247588763155, 2783.02282988684, 0.003280182295545751, 0.003280182295545751),
    @("sCs", "Thbs2", "Itgb1", "sCs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 424.9046499105854, 3824.141849195268, 0.004507286916469871, 0.004507286916469871),
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value2 = $val
        $colIndex++
    }
    $rowIndex++
}
